$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper format sources already present in the sheet (same visual "table"
# pattern used throughout the document):
#   206            -> section title row (merged, bold 14pt on fill)
#   208            -> column header row (Field Name/Data type/Length/.../Comment)
#   209            -> "ID" row style (s=3 across all 5 cells)
#   210            -> normal field row style (s=5,5,5,3,3)
#   58             -> field row style with D also s=5 (s=5,5,5,5,3)
# ---------------------------------------------------------------------------

function Copy-RowFormat($srcRange, $dstRange) {
    $ws.Range($srcRange).Copy() | Out-Null
    $ws.Range($dstRange).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# =============================== Table :-Product Review ====================
Copy-RowFormat "A206:E206" "A216:E216"
$ws.Rows(216).RowHeight = 18.75
$ws.Range("A216").Value = "Table :-Product Review"
$ws.Range("A216:E216").Merge() | Out-Null

Copy-RowFormat "A208:E208" "A218:E218"
$ws.Range("A218").Value = "Field Name"
$ws.Range("B218").Value = "Data type"
$ws.Range("C218").Value = "Length"
$ws.Range("D218").Value = "Nullable"
$ws.Range("E218").Value = "Comment"

Copy-RowFormat "A209:E209" "A219:E219"
$ws.Range("A219").Value = "ID"
$ws.Range("B219").Value = "int"
$ws.Range("C219").Value = 10
$ws.Range("D219").Value = "NO"
$ws.Range("E219").Value = "Primary key,Auto Increment"

Copy-RowFormat "A210:E210" "A220:E220"
$ws.Range("A220").Value = "Product_ID"
$ws.Range("B220").Value = "int"
$ws.Range("C220").Value = 10
$ws.Range("D220").Value = "NO"
$ws.Range("E220").Value = "Reruired"

Copy-RowFormat "A58:E58" "A221:E221"
$ws.Range("A221").Value = "Product_Rate"
$ws.Range("B221").Value = "int"
$ws.Range("C221").Value = 10
$ws.Range("D221").Value = "NO"
$ws.Range("E221").Value = "Reruired"

Copy-RowFormat "A210:E210" "A222:E222"
$ws.Range("A222").Value = "Product_Review"
$ws.Range("B222").Value = "varchar"
$ws.Range("C222").Value = 255
$ws.Range("D222").Value = "NO"
$ws.Range("E222").Value = "Reruired"

Copy-RowFormat "A212:E212" "A223:E223"
$ws.Range("A223").Value = "Created_At"
$ws.Range("B223").Value = "daretime"
$ws.Range("C223").Value = $null
$ws.Range("D223").Value = "YES"
$ws.Range("E223").Value = "Date & Time When  Created"

Copy-RowFormat "A213:E213" "A224:E224"
$ws.Range("A224").Value = "Modified_At"
$ws.Range("B224").Value = "datetime"
$ws.Range("C224").Value = $null
$ws.Range("D224").Value = "YES"
$ws.Range("E224").Value = "Date & Time When  Updated"

# =============================== Table :-All NewsLetter =====================
Copy-RowFormat "A206:E206" "A227:E227"
$ws.Rows(227).RowHeight = 18.75
$ws.Range("A227").Value = "Table :-All NewsLetter"
$ws.Range("A227:E227").Merge() | Out-Null

Copy-RowFormat "A208:E208" "A229:E229"
$ws.Range("A229").Value = "Field Name"
$ws.Range("B229").Value = "Data type"
$ws.Range("C229").Value = "Length"
$ws.Range("D229").Value = "Nullable"
$ws.Range("E229").Value = "Comment"

Copy-RowFormat "A209:E209" "A230:E230"
$ws.Range("A230").Value = "ID"
$ws.Range("B230").Value = "int"
$ws.Range("C230").Value = 10
$ws.Range("D230").Value = "NO"
$ws.Range("E230").Value = "Primary key,Auto Increment"

Copy-RowFormat "A210:E210" "A231:E231"
$ws.Range("A231").Value = "Title"
$ws.Range("B231").Value = "varchar"
$ws.Range("C231").Value = 50
$ws.Range("D231").Value = "NO"
$ws.Range("E231").Value = "Reruired"

Copy-RowFormat "A58:E58" "A232:E232"
$ws.Range("A232").Value = "Description"
$ws.Range("B232").Value = "varchar"
$ws.Range("C232").Value = 255
$ws.Range("D232").Value = "NO"
$ws.Range("E232").Value = "Reruired"

Copy-RowFormat "A212:E212" "A233:E233"
$ws.Range("A233").Value = "Created_At"
$ws.Range("B233").Value = "daretime"
$ws.Range("C233").Value = $null
$ws.Range("D233").Value = "YES"
$ws.Range("E233").Value = "Date & Time When  Created"

Copy-RowFormat "A213:E213" "A234:E234"
$ws.Range("A234").Value = "Modified_At"
$ws.Range("B234").Value = "datetime"
$ws.Range("C234").Value = $null
$ws.Range("D234").Value = "YES"
$ws.Range("E234").Value = "Date & Time When  Updated"

# =============================== Table :- NewsLetter User ===================
Copy-RowFormat "A206:E206" "A237:E237"
$ws.Rows(237).RowHeight = 18.75
$ws.Range("A237").Value = "Table :- NewsLetter User"
$ws.Range("A237:E237").Merge() | Out-Null

Copy-RowFormat "A208:E208" "A239:E239"
$ws.Range("A239").Value = "Field Name"
$ws.Range("B239").Value = "Data type"
$ws.Range("C239").Value = "Length"
$ws.Range("D239").Value = "Nullable"
$ws.Range("E239").Value = "Comment"

Copy-RowFormat "A209:E209" "A240:E240"
$ws.Range("A240").Value = "ID"
$ws.Range("B240").Value = "int"
$ws.Range("C240").Value = 10
$ws.Range("D240").Value = "NO"
$ws.Range("E240").Value = "Primary key,Auto Increment"

Copy-RowFormat "A210:E210" "A241:E241"
$ws.Range("A241").Value = "Email"
$ws.Range("B241").Value = "varchar"
$ws.Range("C241").Value = 50
$ws.Range("D241").Value = "NO"
$ws.Range("E241").Value = "Reruired"

Copy-RowFormat "A212:E212" "A242:E242"
$ws.Range("A242").Value = "Created_At"
$ws.Range("B242").Value = "daretime"
$ws.Range("C242").Value = $null
$ws.Range("D242").Value = "YES"
$ws.Range("E242").Value = "Date & Time When  Created"

# ---------------------------------------------------------------------------
# Update the view so the freshly-added tables are in frame (mirrors the
# commit's sheetView change: topLeftCell="A226", activeCell/sqref="A234").
# ---------------------------------------------------------------------------
$ws.Range("A234").Select()

Write-Output "applied global search bar / product review tables"
